$d = $word.ActiveDocument

# Anchor: last existing paragraph ("Now, the scenario...")
$anchor = $d.Paragraphs.Last

# Insert 9 new paragraphs of text after the anchor; each inherits the
# ListParagraph style / numPr formatting automatically from the paragraph
# mark it was split from.
$texts = @(
    'What makes this scenario more interesting is that now there are two ways in which an estimator can obtain better estimates. Either by removing added pseudorandom noises from measurements that it has access to or by simply fusing more measurements. Of course, when we want to guarantee some privilege gained from holding keys, we want the first case to always make estimates better, but for fusing to only better estimates when corresponding keys are held as well.',
    'To formalise the problem similarly to the single sensor case, we consider linear systems to aid cryptographic proofs about differences in estimations, but now have n different measurements and modification, for now, can be considered an action that produces some modified measurements using a secret key unique to each sensor. A small example on the right shows the measurements of 3 noises sensors of some path that estimators which to estimate.',
    'In addition to the standard linear model assumption we formalise some notation and add assumptions on which sensors can be accessed to simplify notation as well as support number generation which I will get to later. We define an estimator as having access to q measurements and access to a subset p of the measurement keys.',
    'The added assumptions we make are on this subset. We say that access to measurements is sequential, that is when q equals 3, an estimator has access to measurements 1,2 and 3, and similarly for privilege, where privilege encompasses having keys up to and including the privilege value p. Lastly, we also assume that when an estimator holds a key it also has access to the measurements associated with that key so simplify computing bounds.',
    'This in essence lets us group estimators into three classes. Those with no keys that are unprivileged, those that only have access to the measurements for which they hold keys and those with more measurements than keys, which are all still sequential. It can also be noted here that the only assumption that affects the methods presented is that keys that are held are sequential, while the other assumptions are just there to make notation and explanation easier.',
    'Now, with the problem defined, we can describe the performances that we are interested in. That is, how to capture this gain in performance when keys are held or when additional measurements are fused from a cryptographic point of view.',
    'The first difference we’re interested in the performance loss lower bound. It is the difference between an estimator that holds keys for its measurements and one that holds no keys and is specific to a privilege p. Note that if the unprivileged estimator has access to fewer measurements or the privileged one to more measurements for which is does not hold a key, the bound remains a lower bound [explain on pic].',
    'The second bound is the performance gain upper bound. It is the difference between the same estimator that hold keys to its measurements and one that has the same number of keys but also access to the remaining measurements. This is an upper bound as access to fewer remaining measurements can only decrease this difference.',
    'The goal is for these bounds to be computable and capture minimum decrease in performance possible when not being privileged and the maximum performance increase when fusing unprivileged measurements. And needs to be computable for each privilege p you may be interested in.'
)

$firstNewStart = -1
foreach ($t in $texts) {
    $tail = $d.Paragraphs.Last
    $tail.Range.InsertParagraphAfter()
    $p = $d.Paragraphs.Last
    if ($firstNewStart -eq -1) {
        $firstNewStart = $p.Range.Start
    }
    $p.Range.Text = $t
}

# The _GoBack bookmark (Word auto-tracks the last edited location) ends up
# collapsed right before the text of the first newly typed paragraph
# ("What makes this scenario..."), matching where the existing bookmark
# used to sit (end of the previously-last paragraph) before new text push
# beyond it.
$collapsed = $d.Range($firstNewStart, $firstNewStart)
$d.Bookmarks.Add('_GoBack', $collapsed)

# Trailing empty paragraph left at the end (new bullet point, no text yet)
$tail = $d.Paragraphs.Last
$tail.Range.InsertParagraphAfter()
$trailing = $d.Paragraphs.Last
$trailing.SpaceBefore = 0
$trailing.SpaceAfter = 8

Write-Output ('Paragraphs: ' + $d.Paragraphs.Count)

